$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-03-30 Sunday" "2025-03-31 Monday"
Replace-Text "833×4=" "893×2="
Replace-Text "763×9=" "826×9="
Replace-Text "673×8=" "799×7="
Replace-Text "275×5=" "995×6="
Replace-Text "619×7=" "967×8="
Replace-Text "916×7=" "285×8="
Replace-Text "266×2=" "383×2="
Replace-Text "899×5=" "824×7="
Replace-Text "885×9=" "292×2="
Replace-Text "790×3=" "549×6="
Replace-Text "467×4=" "741×3="
Replace-Text "422×3=" "722×2="
Replace-Text "432×7=" "466×9="
Replace-Text "592×5=" "919×8="
Replace-Text "366×4=" "549×6="
Replace-Text "543×4=" "548×9="
Replace-Text "458×3=" "109×5="
Replace-Text "515×6=" "149×5="
Replace-Text "968×8=" "396×3="
Replace-Text "849×8=" "820×3="
Replace-Text "304×3=" "333×6="
Replace-Text "472×8=" "980×2="
Replace-Text "430×5=" "823×9="
Replace-Text "721×2=" "101×7="
Replace-Text "478×6=" "557×2="
